$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) and volume (E) columns are treated as text so numeric-looking
# strings like "27.899.71" or "1.000" are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '27.913.25'
$ws.Range("E2").Value = '  +4.95%  '

# Row 3
$ws.Range("D3").Value = '1.777.73'
$ws.Range("E3").Value = '  +3.48%  '

# Row 4
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.18%  '

# Row 5
$ws.Range("D5").Value = '243.30'
$ws.Range("E5").Value = '  +1.17%  '

# Row 6
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.14%  '

# Row 7
$ws.Range("D7").Value = '0.4893'
$ws.Range("E7").Value = '  -0.53%  '

# Row 8
$ws.Range("D8").Value = '0.2647'
$ws.Range("E8").Value = '  +1.93%  '

# Row 9
$ws.Range("D9").Value = '0.06234'
$ws.Range("E9").Value = '  +0.47%  '

# Row 10
$ws.Range("D10").Value = '1.785.48'
$ws.Range("E10").Value = '  +3.33%  '

# Row 11
$ws.Range("E11").Value = '  +3.39%  '

# Row 12
$ws.Range("D12").Value = '0.07002'
$ws.Range("E12").Value = '  +0.05%  '

# Row 13
$ws.Range("D13").Value = '0.6151'
$ws.Range("E13").Value = '  +1.30%  '

# Row 14
$ws.Range("D14").Value = '4.601'
$ws.Range("E14").Value = '  +2.75%  '

# Row 15
$ws.Range("D15").Value = '79.30'
$ws.Range("E15").Value = '  +3.38%  '

# Row 16
$ws.Range("D16").Value = '27.890.22'
$ws.Range("E16").Value = '  +5.47%  '

# Row 17
$ws.Range("E17").Value = '  +0.20%  '

# Row 18
$ws.Range("D18").Value = '0.9997'
$ws.Range("E18").Value = '  +0.14%  '

# Row 19
$ws.Range("D19").Value = '0.000007193'
$ws.Range("E19").Value = '  +0.70%  '

# Row 20
$ws.Range("E20").Value = '  +3.64%  '

# Row 21
$ws.Range("D21").Value = '2.007.35'
$ws.Range("E21").Value = '  +2.89%  '

# Row 22
$ws.Range("D22").Value = '4.555'
$ws.Range("E22").Value = '  +3.27%  '

# Row 23
$ws.Range("D23").Value = '8.620'
$ws.Range("E23").Value = '  +1.47%  '

# Row 24
$ws.Range("D24").Value = '5.192'
$ws.Range("E24").Value = '  +2.18%  '

# Row 25
$ws.Range("D25").Value = '141.80'
$ws.Range("E25").Value = '  +3.14%  '

# Row 26
$ws.Range("E26").Value = '  +2.18%  '

# Row 27
$ws.Range("D27").Value = '1.856'
$ws.Range("E27").Value = '  +6.60%  '

# Row 28
$ws.Range("D28").Value = '108.97'
$ws.Range("E28").Value = '  +3.28%  '

# Row 29
$ws.Range("E29").Value = '  -0.64%  '

# Row 30
$ws.Range("D30").Value = '4.104'
$ws.Range("E30").Value = '  +4.87%  '

# Row 31
$ws.Range("D31").Value = '0.08227'
$ws.Range("E31").Value = '  +3.67%  '

# Row 32
$ws.Range("D32").Value = '3.775'
$ws.Range("E32").Value = '  +3.74%  '

# Row 33
$ws.Range("D33").Value = '0.04738'
$ws.Range("E33").Value = '  +5.17%  '

# Row 34
$ws.Range("D34").Value = '1.055'
$ws.Range("E34").Value = '  +5.75%  '

# Row 35
$ws.Range("E35").Value = '  -0.71%  '

# Row 36
$ws.Range("D36").Value = '0.6405'
$ws.Range("E36").Value = '  +2.44%  '

# Row 37
$ws.Range("D37").Value = '0.9415'
$ws.Range("E37").Value = '  +0.64%  '

# Row 38
$ws.Range("D38").Value = '2.586'
$ws.Range("E38").Value = '  +7.09%  '

# Row 39
$ws.Range("D39").Value = '2.042'
$ws.Range("E39").Value = '  +1.73%  '

# Row 40
$ws.Range("D40").Value = '5.889'
$ws.Range("E40").Value = '  +6.89%  '

# Row 41
$ws.Range("D41").Value = '0.01534'
$ws.Range("E41").Value = '  +2.06%  '

# Row 42
$ws.Range("E42").Value = '  +0.24%  '

# Row 43
$ws.Range("D43").Value = '100.15'
$ws.Range("E43").Value = '  +0.74%  '

# Row 44
$ws.Range("D44").Value = '0.3933'
$ws.Range("E44").Value = '  +2.48%  '

# Row 45
$ws.Range("D45").Value = '7.162'
$ws.Range("E45").Value = '  +3.24%  '

# Row 46
$ws.Range("D46").Value = '0.1191'
$ws.Range("E46").Value = '  +3.27%  '

# Row 47
$ws.Range("D47").Value = '0.05408'
$ws.Range("E47").Value = '  +0.72%  '

# Row 48
$ws.Range("D48").Value = '7.981'
$ws.Range("E48").Value = '  +2.58%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.279'
$ws.Range("E49").Value = '  +4.83%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '30.35'
$ws.Range("E50").Value = '  +1.01%  '

# Row 51
$ws.Range("D51").Value = '52.40'
$ws.Range("E51").Value = '  +1.92%  '
